$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6994775
$ws.Range("I33").Value = 2071.7273
$ws.Range("K33").Value = 2071.7273
$ws.Range("M33").Value = -1842.7273

$ws.Range("H57").Value = 35590
$ws.Range("J57").Value = 35590
$ws.Range("L57").Value = 106770
$ws.Range("N57").Value = -107768

$ws.Range("H116").Value = 7474.55
$ws.Range("I116").Value = 11700.5
$ws.Range("J116").Value = 3248.6
$ws.Range("K116").Value = 11700.5
$ws.Range("L116").Value = 3248.6
$ws.Range("M116").Value = -8258.5
$ws.Range("N116").Value = -10132.6

$ws.Range("H129").Value = 1035.3043
$ws.Range("J129").Value = 1106.8572
$ws.Range("L129").Value = 3320.5716
$ws.Range("N129").Value = -13320.5716

$ws.Range("H132").Value = 2231.7083
$ws.Range("I132").Value = 1846.0714
$ws.Range("J132").Value = 2771.6
$ws.Range("K132").Value = 5538.2142
$ws.Range("L132").Value = 8314.799999999999
$ws.Range("M132").Value = -3008.2142
$ws.Range("N132").Value = -13374.8

$ws.Range("H138").Value = 1381.12
$ws.Range("I138").Value = 837.3022999999999
$ws.Range("J138").Value = 1791.3684
$ws.Range("K138").Value = 2511.9069
$ws.Range("L138").Value = 5374.1052
$ws.Range("M138").Value = 2628.0931
$ws.Range("N138").Value = -15654.1052

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1043.37
$ws.Range("I32").Value = 1010.41113
$ws.Range("J32").Value = 1340
$ws.Range("K32").Value = 1010.41113
$ws.Range("L32").Value = 1340
$ws.Range("M32").Value = -723.41113
$ws.Range("N32").Value = -1914

$ws.Range("H61").Value = 3993.7368
$ws.Range("I61").Value = 4432.968
$ws.Range("J61").Value = 2048.5715
$ws.Range("K61").Value = 4432.968
$ws.Range("L61").Value = 2048.5715
$ws.Range("M61").Value = -4220.968
$ws.Range("N61").Value = -2472.5715

$ws.Range("H74").Value = 1432.079
$ws.Range("I74").Value = 1268.3
$ws.Range("J74").Value = 2046.25
$ws.Range("K74").Value = 1268.3
$ws.Range("L74").Value = 2046.25
$ws.Range("M74").Value = -394.3
$ws.Range("N74").Value = -3794.25

$ws.Range("H77").Value = 1432.079
$ws.Range("I77").Value = 1268.3
$ws.Range("J77").Value = 2046.25
$ws.Range("K77").Value = 6341.5
$ws.Range("L77").Value = 10231.25
$ws.Range("M77").Value = -1973.5
$ws.Range("N77").Value = -18967.25

$ws.Range("H132").Value = 3556.4167
$ws.Range("I132").Value = 2335.074
$ws.Range("J132").Value = 7220.4443
$ws.Range("K132").Value = 7005.222
$ws.Range("L132").Value = 21661.3329
$ws.Range("M132").Value = -4475.222
$ws.Range("N132").Value = -26721.3329

$ws.Range("H136").Value = 3993.7368
$ws.Range("I136").Value = 4432.968
$ws.Range("J136").Value = 2048.5715
$ws.Range("K136").Value = 13298.904
$ws.Range("L136").Value = 6145.7145
$ws.Range("M136").Value = -10748.904
$ws.Range("N136").Value = -11245.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4794.8335
$ws.Range("I134").Value = 6693.2856
$ws.Range("J134").Value = 2896.3809
$ws.Range("K134").Value = 20079.8568
$ws.Range("L134").Value = 8689.1427
$ws.Range("M134").Value = -17544.8568
$ws.Range("N134").Value = -13759.1427

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 208949.77
$ws.Range("I31").Value = 1450
$ws.Range("J31").Value = 820528.0600000001
$ws.Range("K31").Value = 1450
$ws.Range("L31").Value = 820528.0600000001
$ws.Range("M31").Value = -1155
$ws.Range("N31").Value = -821118.0600000001

$ws.Range("H34").Value = 208949.77
$ws.Range("I34").Value = 1450
$ws.Range("J34").Value = 820528.0600000001
$ws.Range("K34").Value = 1450
$ws.Range("L34").Value = 820528.0600000001
$ws.Range("M34").Value = -1248
$ws.Range("N34").Value = -820932.0600000001

$ws.Range("H58").Value = 1669.6904
$ws.Range("I58").Value = 1121
$ws.Range("J58").Value = 2122.9565
$ws.Range("K58").Value = 1121
$ws.Range("L58").Value = 2122.9565
$ws.Range("M58").Value = -918
$ws.Range("N58").Value = -2528.9565

$ws.Range("H122").Value = 3242.8
$ws.Range("I122").Value = 2250
$ws.Range("J122").Value = 3904.6667
$ws.Range("K122").Value = 6750
$ws.Range("L122").Value = 11714.0001
$ws.Range("M122").Value = -4300
$ws.Range("N122").Value = -16614.0001

$ws.Range("H132").Value = 2398.0222
$ws.Range("I132").Value = 1896.2812
$ws.Range("J132").Value = 3633.077
$ws.Range("K132").Value = 5688.8436
$ws.Range("L132").Value = 10899.231
$ws.Range("M132").Value = -3158.8436
$ws.Range("N132").Value = -15959.231

$ws.Range("H134").Value = 2670.9744
$ws.Range("I134").Value = 3526.8635
$ws.Range("J134").Value = 1563.3529
$ws.Range("K134").Value = 10580.5905
$ws.Range("L134").Value = 4690.0587
$ws.Range("M134").Value = -8045.5905
$ws.Range("N134").Value = -9760.058700000001

$ws.Range("H136").Value = 1669.6904
$ws.Range("I136").Value = 1121
$ws.Range("J136").Value = 2122.9565
$ws.Range("K136").Value = 3363
$ws.Range("L136").Value = 6368.869499999999
$ws.Range("M136").Value = -813
$ws.Range("N136").Value = -11468.8695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 35714384
$ws.Range("I26").Value = 132
$ws.Range("J26").Value = 45454630
$ws.Range("K26").Value = 396
$ws.Range("L26").Value = 136363890
$ws.Range("M26").Value = -108
$ws.Range("N26").Value = -136364466

$ws.Range("H131").Value = 1887741.9
$ws.Range("J131").Value = 1012.9167
$ws.Range("L131").Value = 3038.7501
$ws.Range("N131").Value = -13118.7501

$ws.Range("H132").Value = 1091902.4
$ws.Range("I132").Value = 2037.375
$ws.Range("J132").Value = 1963794.2
$ws.Range("K132").Value = 18336.375
$ws.Range("L132").Value = 17674147.8
$ws.Range("M132").Value = -15806.375
$ws.Range("N132").Value = -17679207.8

$ws.Range("H134").Value = 9616.359
$ws.Range("I134").Value = 12063.8
$ws.Range("J134").Value = 8772.414000000001
$ws.Range("K134").Value = 36191.39999999999
$ws.Range("L134").Value = 26317.242
$ws.Range("M134").Value = -31121.39999999999
$ws.Range("N134").Value = -36457.242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 70009
$ws.Range("J22").Value = 70009
$ws.Range("L22").Value = 70009
$ws.Range("N22").Value = -71067

$ws.Range("H126").Value = 4039.6863
$ws.Range("J126").Value = 2570.1538
$ws.Range("L126").Value = 7710.4614
$ws.Range("N126").Value = -12650.4614

$ws.Range("H132").Value = 2640.4167
$ws.Range("I132").Value = 2206.5
$ws.Range("J132").Value = 3322.2856
$ws.Range("K132").Value = 6619.5
$ws.Range("L132").Value = 9966.856800000001
$ws.Range("M132").Value = -4089.5
$ws.Range("N132").Value = -15026.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 25855.414
$ws.Range("I7").Value = 44973.348
$ws.Range("J7").Value = 1426.9445
$ws.Range("K7").Value = 44973.348
$ws.Range("L7").Value = 1426.9445
$ws.Range("M7").Value = -44861.348
$ws.Range("N7").Value = -1650.9445

$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 9999
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 9999
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -10385

$ws.Range("H119").Value = 29000
$ws.Range("J119").Value = 29000
$ws.Range("L119").Value = 29000
$ws.Range("N119").Value = -38676

$ws.Range("H122").Value = 4074702.8
$ws.Range("I122").Value = 4765263.5
$ws.Range("J122").Value = 2003020
$ws.Range("K122").Value = 14295790.5
$ws.Range("L122").Value = 6009060
$ws.Range("M122").Value = -14293340.5
$ws.Range("N122").Value = -6013960

$ws.Range("H126").Value = 25855.414
$ws.Range("I126").Value = 44973.348
$ws.Range("J126").Value = 1426.9445
$ws.Range("K126").Value = 134920.044
$ws.Range("L126").Value = 4280.833500000001
$ws.Range("M126").Value = -132450.044
$ws.Range("N126").Value = -9220.833500000001

$ws.Range("H132").Value = 22271260
$ws.Range("I132").Value = 28130728
$ws.Range("K132").Value = 84392184
$ws.Range("M132").Value = -84389654

$ws.Range("H136").Value = 8812.799999999999
$ws.Range("I136").Value = 6289.2144
$ws.Range("J136").Value = 14701.167
$ws.Range("K136").Value = 18867.6432
$ws.Range("L136").Value = 44103.501
$ws.Range("M136").Value = -16317.6432
$ws.Range("N136").Value = -49203.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18819.322
$ws.Range("I132").Value = 22920.422
$ws.Range("J132").Value = 2042.091
$ws.Range("K132").Value = 68761.266
$ws.Range("L132").Value = 6126.272999999999
$ws.Range("M132").Value = -66231.266
$ws.Range("N132").Value = -11186.273

$ws.Range("H136").Value = 8931322
$ws.Range("I136").Value = 3266.4062
$ws.Range("J136").Value = 20835396
$ws.Range("K136").Value = 9799.2186
$ws.Range("L136").Value = 62506188
$ws.Range("M136").Value = -7249.2186
$ws.Range("N136").Value = -62511288
